# Stopped duplicate warrants from being generated.
#
# The {seizableItems} placeholder was a plain (non-looping) merge field,
# which caused only the first seizable item to be rendered and, worse,
# could cause the downstream templating engine to duplicate the whole
# warrant. This converts both occurrences of the {seizableItems}
# placeholder (in the Affidavit and in the Seizure Warrant sections)
# into a proper {#seizableItems} ... {text} ... {/seizableItems} loop,
# matching the pattern already used elsewhere in the document (e.g.
# {#affiantHerosheet} / {#factsAndCircumstances}). It also drops the
# paragraph's direct line-spacing override, which is no longer needed
# once the paragraph is followed by the two new loop-control paragraphs.

$d = $word.ActiveDocument

function Get-PkgXml([string]$bodyInnerXml) {
    return '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyInnerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# The two trailing paragraphs that get appended after every
# {#seizableItems} opening tag are identical in both locations.
$loopBodyAndClose =
    '<w:p>' +
        '<w:pPr>' +
            '<w:pStyle w:val="Body - Keep Lines"/>' +
            '<w:ind w:firstLine="1417"/>' +
            '<w:rPr><w:b w:val="1"/><w:bCs w:val="1"/></w:rPr>' +
        '</w:pPr>' +
        '<w:r>' +
            '<w:rPr><w:b w:val="1"/><w:bCs w:val="1"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr>' +
            '<w:t>{text}</w:t>' +
        '</w:r>' +
    '</w:p>' +
    '<w:p>' +
        '<w:pPr>' +
            '<w:pStyle w:val="Body - Keep Lines"/>' +
            '<w:ind w:firstLine="1417"/>' +
        '</w:pPr>' +
        '<w:r>' +
            '<w:rPr><w:b w:val="1"/><w:bCs w:val="1"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr>' +
            '<w:t>{/seizableItems}</w:t>' +
        '</w:r>' +
    '</w:p>'

# ---------------------------------------------------------------------
# Locate both paragraphs that still hold the bare {seizableItems} token.
# Walking back-to-front means inserting/replacing content earlier in the
# story never invalidates the index of a paragraph we have not visited
# yet.
# ---------------------------------------------------------------------
$targets = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*{seizableItems}*") {
        $targets += $i
    }
}

for ($k = $targets.Count - 1; $k -ge 0; $k--) {
    $idx = $targets[$k]
    $p = $d.Paragraphs.Item($idx)
    $hasTab = $p.Range.Text -like "*`t*"

    if ($hasTab) {
        # First occurrence: paragraph already starts with a tab run
        # before the placeholder text, so only the placeholder text
        # itself needs braces added, and the <w:spacing> override is
        # dropped.
        $firstPara =
            '<w:p>' +
                '<w:pPr>' +
                    '<w:pStyle w:val="Body"/>' +
                    '<w:keepLines w:val="1"/>' +
                    '<w:ind w:firstLine="850"/>' +
                    '<w:rPr><w:b w:val="1"/><w:bCs w:val="1"/></w:rPr>' +
                '</w:pPr>' +
                '<w:r>' +
                    '<w:rPr><w:b w:val="0"/><w:bCs w:val="0"/></w:rPr>' +
                    '<w:tab/>' +
                '</w:r>' +
                '<w:r>' +
                    '<w:rPr><w:b w:val="1"/><w:bCs w:val="1"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr>' +
                    '<w:t>{#seizableItems}</w:t>' +
                '</w:r>' +
            '</w:p>'
    } else {
        # Second occurrence: paragraph has no leading tab run yet, so a
        # new non-bold tab run is introduced ahead of the (now looped)
        # placeholder, and the <w:spacing> override is dropped.
        $firstPara =
            '<w:p>' +
                '<w:pPr>' +
                    '<w:pStyle w:val="Body"/>' +
                    '<w:keepLines w:val="1"/>' +
                    '<w:ind w:firstLine="850"/>' +
                    '<w:rPr><w:b w:val="1"/><w:bCs w:val="1"/></w:rPr>' +
                '</w:pPr>' +
                '<w:r>' +
                    '<w:rPr><w:b w:val="0"/><w:bCs w:val="0"/></w:rPr>' +
                    '<w:tab/>' +
                '</w:r>' +
                '<w:r>' +
                    '<w:rPr><w:b w:val="1"/><w:bCs w:val="1"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr>' +
                    '<w:t>{#seizableItems}</w:t>' +
                '</w:r>' +
            '</w:p>'
    }

    $xml = Get-PkgXml ($firstPara + $loopBodyAndClose)
    $p.Range.InsertXML($xml)
}

Write-Host "Converted" $targets.Count "{seizableItems} placeholder(s) into loops"
